$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Gdf1", "Bmpr1a", "ECs", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 3, 1, 2.535712666666667, 7.607138, 0.04494879354621957, 0.05070282964779482, 0.01742034602, 0.15678311418, 0.009040554254253255, 0.01019786397167053),
    @("FAPs", "Gdf1", "Bmpr1a", "FAPs", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 3, 1, 32.24261766666667, 96.72785300000001, 0.5715421877013505, 0.6447070965264385, 0.22150678337, 1.99356105033, 0.1149543235503199, 0.1296699872101365),
    @("FAPs", "Gdf1", "Bmpr1a", "Inflammatory-Mac", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 3, 1, 1.538811333333333, 4.616434, 0.02727742533206951, 0.03076929413956839, 0.01057163386, 0.09514470474, 0.005486310625386233, 0.006188630463414083),
    @("FAPs", "Gdf1", "Bmpr1a", "MuSCs", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 2, 1, 19.206297, 38.412594, 0.3404565075487166, 0.2560262755732715, 0.13194726039, 0.7916835623399999, 0.06847604317884132, 0.05149458421958529),
    @("FAPs", "Gdf1", "Bmpr1a", "Resolving-Mac", 1, 0.3333333333333333, 0.00687, 0.02061, 0.2011300758263314, 0.2011300758263313, 3, 1, 0.8899256666666666, 2.669777, 0.0157750858716439, 0.01779450411292666, 0.00611378933, 0.05502410396999999, 0.003172844217530626, 0.003579009961524904),
    @("Resolving-Mac", "Gdf1", "Bmpr1a", "ECs", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 3, 1, 2.535712666666667, 7.607138, 0.04494879354621957, 0.05070282964779482, 0.06919199153533334, 0.622727923818, 0.03590823929196632, 0.04050496567612429),
    @("Resolving-Mac", "Gdf1", "Bmpr1a", "FAPs", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 3, 1, 32.24261766666667, 96.72785300000001, 0.5715421877013505, 0.6447070965264385, 0.8798043082703334, 7.918238774433001, 0.4565878641510306, 0.515037109316302),
    @("Resolving-Mac", "Gdf1", "Bmpr1a", "Inflammatory-Mac", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 3, 1, 1.538811333333333, 4.616434, 0.02727742533206951, 0.03076929413956839, 0.04198954485266667, 0.377905903674, 0.02179111470668328, 0.02458066367615431),
    @("Resolving-Mac", "Gdf1", "Bmpr1a", "MuSCs", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 2, 1, 19.206297, 38.412594, 0.3404565075487166, 0.2560262755732715, 0.524082226239, 3.144493357434, 0.2719804643698753, 0.2045316913536862),
    @("Resolving-Mac", "Gdf1", "Bmpr1a", "Resolving-Mac", 1, 0.3333333333333333, 0.027287, 0.081861, 0.7988699241736688, 0.7988699241736686, 3, 1, 0.8899256666666666, 2.669777, 0.0157750858716439, 0.01779450411292666, 0.02428340166633334, 0.218550614997, 0.01260224165411328, 0.01421549415140176)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
